# Berichtsheft KW12 2025: add a new line "Schreibtischtests" right after the
# existing "Bedeutung von künstlicher Intelligenz für Unternehmen" entry in
# the Freitag / "Ausgeführte Arbeiten" table cell.

$d = $word.ActiveDocument

$anchorText = "Bedeutung von künstlicher Intelligenz für Unternehmen"
$newLineText = "Schreibtischtests"

# Locate the table cell that holds the anchor paragraph and search inside it
# so the Find stays scoped to that single cell.
$targetCell = $null
for ($ti = 1; $ti -le $d.Tables.Count -and $targetCell -eq $null; $ti++) {
    $tbl = $d.Tables.Item($ti)
    for ($ri = 1; $ri -le $tbl.Rows.Count -and $targetCell -eq $null; $ri++) {
        for ($ci = 1; $ci -le $tbl.Columns.Count -and $targetCell -eq $null; $ci++) {
            try {
                $cell = $tbl.Cell($ri, $ci)
            } catch {
                $cell = $null
            }
            if ($cell -ne $null -and $cell.Range.Text -like "*$anchorText*") {
                $targetCell = $cell
            }
        }
    }
}

if ($targetCell -eq $null) {
    throw "Could not find the target cell containing '$anchorText'"
}

$rng = $targetCell.Range
$found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0, $false)

if (-not $found) {
    throw "Could not find anchor text '$anchorText' inside target cell"
}

# Collapse to the end of the matched text, then append a new paragraph
# carrying the same run/paragraph formatting (size 20, 260-twip exact
# spacing) that the surrounding paragraphs already use.
$rng.Collapse(0)
$rng.InsertAfter([char]13 + $newLineText)

Write-Output "Inserted new paragraph '$newLineText'."
